$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Lee Harvey Oswald killed President Kennedy on his own."
$ws.Range("B2").Value = "'TRUE"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = 0.97

# Row 3
$ws.Range("A3").Value = "The federal authorities routinely collect data on phone calls, emails and other electronic traffic that Americans generate, regardless of whether they have any bearing on a counterterrorism investigation."
$ws.Range("B3").Value = "'FALSE"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = 0.7

# Row 4
$ws.Range("A4").Value = "North Korea received Soviet and Chinese aid in developing its chemical industry."
$ws.Range("B4").Value = "'TRUE"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = 0.85

# Row 5
$ws.Range("A5").Value = "Despite the severity of the impacts of Hurricane Maria, the Government of Puerto Rico initially reported that only 64 people died in the hurricane. Numerous media outlets harshly criticized the government for suppressing the true death toll."
$ws.Range("B5").Value = "'TRUE"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = 0.95

# Row 6
$ws.Range("A6").Value = "US Representative Anthony Weiner's Twitter account linked to an inappropriate photograph. Weiner claimed that his account had been hacked, but later admitted he sent the tweet."
$ws.Range("B6").Value = "'TRUE"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = 0.99

# Row 7
$ws.Range("A7").Value = "Volkswagen had intentionally programmed diesel engines to activate their emissions controls only during laboratory emissions testing, which caused the vehicles' NOx output to meet US standards during regulatory testing. However, the vehicles emitted up to 40 times more NOx in real-world driving."
$ws.Range("B7").Value = "'TRUE"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = 1

# Row 8
$ws.Range("A8").Value = "The assassination of Julius Caesar was a conspiracy orchestrated by Marcus Brutus and other Roman senators."
$ws.Range("B8").Value = "'TRUE"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 1

# Row 9
$ws.Range("A9").Value = "The Watergate hotel room used by Democratic National Committee was bugged by Republican officials, operating at the behest of the White House."
$ws.Range("B9").Value = "'TRUE"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 0.8

# Row 10
$ws.Range("A10").Value = "In the 1950s, the Central Intelligence Agency administered LSD and related drugs under Project MKULTRA, in an effort to investigate the possibility of “mind control”."
$ws.Range("B10").Value = "'TRUE"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = 0.9
